# Add two new columns "I0" (I) and "IF" (J) to the sheet, mirroring the
# existing header style used by column H, and fill in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ------------------------------------------------
# Copy the formatting (bold font, border, centered/top alignment) that the
# existing header cell H1 uses onto the two new header cells, then set the
# header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-45) ------------------------------------------------
# Each entry is: row, I-value, J-value
$data = @(
    @(2,7,7),
    @(3,5,5),
    @(4,6,7),
    @(5,5,5),
    @(6,6,6),
    @(7,6,6),
    @(8,7,8),
    @(9,8,8),
    @(10,9,9),
    @(11,9,9),
    @(12,6,7),
    @(13,7,7),
    @(14,6,6),
    @(15,8,8),
    @(16,7,7),
    @(17,7,8),
    @(18,7,7),
    @(19,7,7),
    @(20,5,5),
    @(21,8,8),
    @(22,7,7),
    @(23,7,8),
    @(24,7,7),
    @(25,8,8),
    @(26,5,6),
    @(27,9,9),
    @(28,8,8),
    @(29,6,7),
    @(30,7,7),
    @(31,7,8),
    @(32,8,8),
    @(33,7,7),
    @(34,6,6),
    @(35,5,5),
    @(36,8,8),
    @(37,6,7),
    @(38,9,9),
    @(39,6,6),
    @(40,6,6),
    @(41,8,8),
    @(42,8,8),
    @(43,3,3),
    @(44,3,3),
    @(45,7,7)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal   # column I
    $ws.Cells.Item($r, 10).Value = $jVal  # column J
}

Write-Output "I0/IF columns added"
